$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new "topic" header in I1, matching the style of the
# --- existing bold/bordered header cells (copy format from A1). ---
$ws.Range("A1").Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "topic"

# --- Row 2: update quiz data / counts for the existing quiz_prompt_socratic
# --- materials_5_6 row (quizId 33 -> 57). ---
$ws.Cells.Item(2, 4).Value = '{''quizId'': 57, ''courseId'': 11, ''materialIds'': [5, 6], ''section'': None, ''quizTitle'': ''K-Means, K-NN, and Clustering Algorithms Quiz'', ''questions'': [{''id'': 676, ''question'': ''In K-Means clustering, what is a potential consequence of an inappropriate initialization of cluster centroids?'', ''options'': [''The algorithm may converge to a suboptimal clustering.'', ''The algorithm will always find the global optimum regardless of initialization.'', ''The algorithm will not converge at all.'', ''The algorithm will automatically adjust the number of clusters, K, to compensate.''], ''correctAnswer'': 0, ''explanation'': "K-Means is sensitive to initial centroid placement. Poor initialization can lead to convergence at a local optimum, not the global optimum. The algorithm doesn''t adjust K automatically, and it will always converge (though possibly to a poor solution).", ''source'': {''chunkId'': 705, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''Example of how K-Means could converge to the wrong thing\n\nHow might you make K-Means more likely to converge to the optimal?\n\nHow might you choose K? (Why can’t you measure test error?)\n\nOther ways to evaluate clusters (“cluster validation”)\n\nThere are loa''}, ''grade'': None, ''feedback'': None}, {''id'': 677, ''question'': ''Why is directly measuring test error not a suitable method for choosing the optimal number of clusters (K) in K-Means?'', ''options'': [''Test error is only applicable to classification problems, not clustering.'', ''Clustering is an unsupervised learning technique, and test error requires labeled data.'', ''Test error always decreases as K increases, making it uninformative.'', ''Calculating test error is computationally too expensive for large datasets.''], ''correctAnswer'': 1, ''explanation'': ''Clustering is unsupervised, meaning there are no predefined labels to compare against. Test error requires labeled data to evaluate the performance of the clustering. While computational cost can be a factor, the primary reason is the absence of ground truth labels.'', ''source'': {''chunkId'': 705, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''Example of how K-Means could converge to the wrong thing\n\nHow might you make K-Means more likely to converge to the optimal?\n\nHow might you choose K? (Why can’t you measure test error?)\n\nOther ways to evaluate clusters (“cluster validation”)\n\nThere are loa''}, ''grade'': None, ''feedback'': None}, {''id'': 678, ''question'': ''Lemma 1 establishes a relationship between cost functions in clustering. According to Lemma 1, how does the cost of assigning points in cluster C to an arbitrary point z relate to the cost of assigning them to the mean of C?'', ''options'': [''cost(C, z) is always less than cost(C, mean(C)).'', ''cost(C, z) is always equal to cost(C, mean(C)).'', ''cost(C, z) is equal to cost(C, mean(C)) plus a term proportional to the squared distance between z and the mean of C.'', ''cost(C, z) is equal to cost(C, mean(C)) minus a term proportional to the squared distance between z and the mean of C.''], ''correctAnswer'': 2, ''explanation'': ''Lemma 1 states: cost(C; z) = cost(C, mean(C)) + |C| * ||z - mean(C)||^2. This shows that the cost of assigning points to an arbitrary point z is the cost of assigning them to the mean of C, plus a term that increases with the squared distance between z and the mean.'', ''source'': {''chunkId'': 340, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''To prove Lemma 1, pick a specific choice for X, namely X is a uniform random draw from the points xi in set C. So X has a discrete distribution. What will happen with this choice of X is that the expectation will reduce to the cost we already defined above.\n\nX\n\nEX∥X −z∥2\n2 =\n(prob. that point i is chosen)∥xi −z∥2 2 i:xi C\n\n{ ∈} X 1\n=\n\n{i:xi∈C} |C ∥x | i −z∥2 1\n2 =\ncost(C, z) |C| (2)\n\nand if we use Lemma 2 substituting z to be x¯ (a.k.a., EXX, or mean(C)) and simplify as in (2):\n\n1 EX∥X −x¯∥2\n2 =\ncost(C, mean(C)). C | | (3)\n\nWe had already defined cost earlier, and the choice of X was nice because its expectation is just the cost. Let’s recopy Lemma 2’s statement here, using the x¯ notation:\n\nEX∥X −z∥2\n2 = EX∥X −x¯∥2\n2 + ∥z −x¯∥2\n2.\n\nPlugging in (2) and (3),\n\n1 1\ncost(C, z) =\ncost(C, mean(C)) +\n|C| |C ∥z −x¯∥2 | 2. (4)\n\nMultiplying through,\n\ncost(C; z) = cost(C, mean(C)) + |C| · ∥z −mean(C)∥2\n2.\n\nAnd that’s the statement of Lemma 1. ■\n\nTo really minimize the cost (1), you’d need to try ''}, ''grade'': None, ''feedback'': None}, {''id'': 679, ''question'': ''In the context of Lemma 2, what does the term EX||X - EXX||^2 represent?'', ''options'': [''The expected value of the squared distance between a random variable X and an arbitrary point z.'', ''The variance of the random variable X.'', ''The squared distance between the expected value of X and an arbitrary point z.'', ''The covariance of the random variable X.''], ''correctAnswer'': 1, ''explanation'': ''EX||X - EXX||^2 represents the expected value of the squared difference between X and its expected value (EXX), which is the definition of variance. It measures the spread or dispersion of the random variable X around its mean.'', ''source'': {''chunkId'': 339, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''s one’s pretty much identical to one of the ones we did before).\n\nLemma 2. Let X Rn be any random variable. For any z Rn, we have: ∈ ∈\n\nEX∥X −z∥2\n2 = EX∥X −EXX∥2\n2 + ∥z −EXX∥2\n2.\n\nProof. Let x¯ := EXX.\n\nj −z(j))2 X\n\nEX∥X −z∥2\n2 = EX\n(X(j)\n\nj\n−x¯(j) + x¯(j) −zj)2\nX\n\n= EX\n(X(j)\n\nj\n−x¯(j))2 + EX\n(x¯(j)\n\nj −zj)2 X X\n\n= EX\n(X(j)\n\nj −x¯(j))(x¯(j) −z(j)) X\n\n+2EX\n(X(j)\n\n= EX∥X −x¯∥2\n2 + EX∥x¯ −z∥2\n2 + 0. ■\n\n2''}, ''grade'': None, ''feedback'': None}, {''id'': 680, ''question'': ''Why is an exhaustive search of all possible assignments of data points to clusters generally infeasible for minimizing the cost function in K-Means?'', ''options'': [''The cost function is non-differentiable.'', ''The number of possible assignments grows exponentially with the number of data points and clusters.'', ''The cost function has many local minima.'', ''Exhaustive search requires calculating the inverse of a large matrix.''], ''correctAnswer'': 1, ''explanation'': ''The number of distinct assignments of m data points to K clusters is given by a formula that grows extremely rapidly. S(10, 4) = 34K, S(19, 4) ≈ 10^10. This exponential growth makes exhaustive search computationally impractical.'', ''source'': {''chunkId'': 341, ''fileName'': ''Clustering.pdf'', ''chunkText'': '')) +\n|C| |C ∥z −x¯∥2 | 2. (4)\n\nMultiplying through,\n\ncost(C; z) = cost(C, mean(C)) + |C| · ∥z −mean(C)∥2\n2.\n\nAnd that’s the statement of Lemma 1. ■\n\nTo really minimize the cost (1), you’d need to try all possible assignments of the m data points to K clusters. Uck! The number of distinct assignments is (Jain and Dubes 1988):\n\nK 1\nS(m, K) =\n(\nK! k=1\n−1)K−k K km k\n\nX\n\nS(10, 4) = 34K, S(19, 4)\n1010, ... so not doable. ≈\n\nLet’s try some heuristic gradient-descent-ish method instead.\n\n3''}, ''grade'': None, ''feedback'': None}, {''id'': 681, ''question'': "What is a ''heuristic gradient-descent-ish method'' in the context of K-Means, and why is it used?", ''options'': [''A method that guarantees finding the global minimum of the cost function by exhaustively searching the solution space.'', ''An approximation algorithm that iteratively improves the clustering by making small adjustments, used because finding the absolute minimum is computationally infeasible.'', ''A method for visualizing high-dimensional data in two dimensions.'', ''A technique for preprocessing data to improve the performance of K-Means.''], ''correctAnswer'': 1, ''explanation'': "Because finding the absolute minimum cost is computationally infeasible, K-Means uses a heuristic approach similar to gradient descent. It iteratively adjusts the cluster assignments and centroids to reduce the cost function, but it doesn''t guarantee finding the global minimum.", ''source'': {''chunkId'': 341, ''fileName'': ''Clustering.pdf'', ''chunkText'': '')) +\n|C| |C ∥z −x¯∥2 | 2. (4)\n\nMultiplying through,\n\ncost(C; z) = cost(C, mean(C)) + |C| · ∥z −mean(C)∥2\n2.\n\nAnd that’s the statement of Lemma 1. ■\n\nTo really minimize the cost (1), you’d need to try all possible assignments of the m data points to K clusters. Uck! The number of distinct assignments is (Jain and Dubes 1988):\n\nK 1\nS(m, K) =\n(\nK! k=1\n−1)K−k K km k\n\nX\n\nS(10, 4) = 34K, S(19, 4)\n1010, ... so not doable. ≈\n\nLet’s try some heuristic gradient-descent-ish method instead.\n\n3''}, ''grade'': None, ''feedback'': None}, {''id'': 682, ''question'': ''What is a key advantage of K-Nearest Neighbors (K-NN) regarding model training?'', ''options'': [''K-NN requires extensive training to build a complex model.'', "K-NN involves no explicit training phase; it''s a ''lazy'' learner.", ''K-NN automatically selects the optimal features during training.'', ''K-NN uses a sophisticated training algorithm to minimize overfitting.''], ''correctAnswer'': 1, ''explanation'': "K-NN is a ''lazy'' learner because it doesn''t have an explicit training phase. New training examples can be added easily without retraining the model. The classification or regression is performed at the time of prediction.", ''source'': {''chunkId'': 353, ''fileName'': ''k-nearest neighbpurs.pdf'', ''chunkText'': ''Pros/Cons to K-NN\n\nPros:\n\n• Simple and powerful. No need for tuning complex parameters to build a model.\n\n• No training involved (“lazy”). New training examples can be added easily.\n\nPros/Cons to K-NN\n\nCons:\n\n• Expensive and slow: O(md), m= # examples, d= # dimensions\n\n– To determine the nearest neighbor of a new point x, must compute the distance to all m training examples. Runtime performance is slow, but can be improved.\n\n• Pre-sort training examples into fast data structures\n\n• Compute only an approximate distance\n\n• Remove redundant data (condensing)\n\n4''}, ''grade'': None, ''feedback'': None}, {''id'': 683, ''question'': ''What is a significant computational drawback of using K-Nearest Neighbors (K-NN) for prediction?'', ''options'': [''Calculating distances to all training examples for each prediction can be computationally expensive.'', ''K-NN requires inverting a large covariance matrix, which is computationally intensive.'', ''The training phase of K-NN involves solving a complex optimization problem.'', ''K-NN is prone to overfitting, requiring computationally expensive regularization techniques.''], ''correctAnswer'': 0, ''explanation'': ''To find the nearest neighbors of a new point, K-NN must compute the distance to all m training examples, resulting in a runtime complexity of O(md), where m is the number of examples and d is the number of dimensions. This can be slow for large datasets.'', ''source'': {''chunkId'': 353, ''fileName'': ''k-nearest neighbpurs.pdf'', ''chunkText'': ''Pros/Cons to K-NN\n\nPros:\n\n• Simple and powerful. No need for tuning complex parameters to build a model.\n\n• No training involved (“lazy”). New training examples can be added easily.\n\nPros/Cons to K-NN\n\nCons:\n\n• Expensive and slow: O(md), m= # examples, d= # dimensions\n\n– To determine the nearest neighbor of a new point x, must compute the distance to all m training examples. Runtime performance is slow, but can be improved.\n\n• Pre-sort training examples into fast data structures\n\n• Compute only an approximate distance\n\n• Remove redundant data (condensing)\n\n4''}, ''grade'': None, ''feedback'': None}, {''id'': 684, ''question'': ''How can the computational cost of K-NN be reduced?'', ''options'': [''By increasing the number of neighbors (K) used for prediction.'', ''By pre-sorting training examples into fast data structures or using approximate distance calculations.'', ''By using a more complex distance metric, such as Mahalanobis distance.'', ''By reducing the dimensionality of the data using Principal Component Analysis (PCA).''], ''correctAnswer'': 1, ''explanation'': ''The computational cost of K-NN can be reduced by pre-sorting training examples into fast data structures (e.g., KD-trees) or by computing only an approximate distance. These techniques reduce the number of distance calculations required.'', ''source'': {''chunkId'': 353, ''fileName'': ''k-nearest neighbpurs.pdf'', ''chunkText'': ''Pros/Cons to K-NN\n\nPros:\n\n• Simple and powerful. No need for tuning complex parameters to build a model.\n\n• No training involved (“lazy”). New training examples can be added easily.\n\nPros/Cons to K-NN\n\nCons:\n\n• Expensive and slow: O(md), m= # examples, d= # dimensions\n\n– To determine the nearest neighbor of a new point x, must compute the distance to all m training examples. Runtime performance is slow, but can be improved.\n\n• Pre-sort training examples into fast data structures\n\n• Compute only an approximate distance\n\n• Remove redundant data (condensing)\n\n4''}, ''grade'': None, ''feedback'': None}, {''id'': 685, ''question'': ''In K-NN, what is the impact of choosing a very small value for K?'', ''options'': [''The model will generalize well to unseen data.'', ''The model will be more robust to noise in the data.'', ''The model will overfit the training data and be sensitive to noise.'', ''The model will underfit the training data and fail to capture important patterns.''], ''correctAnswer'': 2, ''explanation'': ''A small value of K makes the model sensitive to noise in the training data because the prediction is based on very few neighbors. This leads to overfitting, where the model performs well on the training data but poorly on unseen data.'', ''source'': {''chunkId'': 49, ''fileName'': ''k-nearest neighbpurs.pdf'', ''chunkText'': ''K-Nearest Neighbors\n\n• For regression: the value for the test eXample becomes the (weighted) average of the values of the K neighbors.\n\nMaking K-NN More Powerful\n\n• A good value for K can be determined by considering a range of K values.\n\n– K too small: we’ll model the noise\n\n– K too large: neighbors include too many points from other classes\n\n• There are problems when there is a spread of distances among the K\xad\nNN. Use a distance-based voting scheme, where closer neighbors have\nmore influence.\n\n• The distance measure has to be meaningful – attributes should be scaled\n\n– Eg. Income varies 10,000-1,000,000 while height varies 1.5-1.8 meters\n\n3''}, ''grade'': None, ''feedback'': None}, {''id'': 686, ''question'': ''In K-NN, why is it important to scale the attributes before calculating distances?'', ''options'': [''Scaling ensures that all attributes have the same mean and standard deviation.'', ''Scaling prevents attributes with larger ranges from dominating the distance calculation.'', ''Scaling reduces the dimensionality of the data.'', ''Scaling makes the data easier to visualize.''], ''correctAnswer'': 1, ''explanation'': ''If attributes have significantly different ranges (e.g., income vs. height), the attribute with the larger range will dominate the distance calculation. Scaling ensures that all attributes contribute equally to the distance, leading to more meaningful neighbor selection.'', ''source'': {''chunkId'': 49, ''fileName'': ''k-nearest neighbpurs.pdf'', ''chunkText'': ''K-Nearest Neighbors\n\n• For regression: the value for the test eXample becomes the (weighted) average of the values of the K neighbors.\n\nMaking K-NN More Powerful\n\n• A good value for K can be determined by considering a range of K values.\n\n– K too small: we’ll model the noise\n\n– K too large: neighbors include too many points from other classes\n\n• There are problems when there is a spread of distances among the K\xad\nNN. Use a distance-based voting scheme, where closer neighbors have\nmore influence.\n\n• The distance measure has to be meaningful – attributes should be scaled\n\n– Eg. Income varies 10,000-1,000,000 while height varies 1.5-1.8 meters\n\n3''}, ''grade'': None, ''feedback'': None}, {''id'': 687, ''question'': ''What is a major issue with K-means clustering when the value of K is changed?'', ''options'': [''The algorithm becomes computationally more expensive.'', ''Cluster membership can change arbitrarily.'', ''The algorithm always converges to the same solution.'', ''The algorithm is guaranteed to find the optimal clustering.''], ''correctAnswer'': 1, ''explanation'': ''A significant issue with K-means is that as K changes, the cluster membership of data points can change arbitrarily. This means that the clusters formed with one value of K may not be related to the clusters formed with a different value of K.'', ''source'': {''chunkId'': 718, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''\ncell Precursors." Journal of Experimental Medicine 203 (2006): 675-87.\n\nwhere each row is a gene, and the columns are different immune cell types.\n\nA major issue with K-means: as K changes, cluster membership can change\narbitrarily. A solution is Hierarch''}, ''grade'': None, ''feedback'': None}, {''id'': 688, ''question'': ''What is the key characteristic of Hierarchical Clustering that addresses the instability of cluster membership in K-Means when K changes?'', ''options'': [''Hierarchical Clustering always produces the same clusters as K-Means.'', ''Hierarchical Clustering builds a hierarchy of clusters, allowing for exploration at different levels of granularity.'', ''Hierarchical Clustering uses a fixed value of K throughout the process.'', ''Hierarchical Clustering is less computationally expensive than K-Means.''], ''correctAnswer'': 1, ''explanation'': ''Hierarchical Clustering creates a hierarchy of clusters by merging clusters at successively higher levels. This allows you to explore the data at different levels of granularity and choose a suitable level of clustering without the arbitrary changes in membership seen in K-Means when K is changed.'', ''source'': {''chunkId'': 719, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''hip can change\narbitrarily. A solution is Hierarchical Clustering.\n\n• clusters at the next level of the hierarchy are created by merging clusters at the next lowest level.\n\n– lowest level: each cluster has 1 example\n\n– highest level: there’s only 1 cluster''}, ''grade'': None, ''feedback'': None}, {''id'': 689, ''question'': ''What are the general characteristics of good cluster validation metrics?'', ''options'': [''They favor small intracluster distances and small intercluster distances.'', ''They favor large intracluster distances and large intercluster distances.'', ''They favor small intracluster distances and large intercluster distances.'', ''They are independent of the distance metric used for clustering.''], ''correctAnswer'': 2, ''explanation'': ''Good cluster validation metrics aim to minimize the distance between points within the same cluster (intracluster distance) and maximize the distance between different clusters (intercluster distance). This indicates that the clusters are well-separated and internally cohesive.'', ''source'': {''chunkId'': 707, ''fileName'': ''Clustering.pdf'', ''chunkText'': ''istance) to the centroid (want it to be small), and intercluster distances between centroids (want it to be large).\n\n• Dunn Index - looks pairwise at minimal intercluster distance (want it to be\nlarge) and maximal intracluster distance (want it to be small''}, ''grade'': None, ''feedback'': None}, {''id'': 690, ''question'': ''According to Lemma 3, what happens to the cost function during the iterations of the K-Means algorithm?'', ''options'': [''The cost monotonically increases.'', ''The cost remains constant.'', ''The cost monotonically decreases.'', ''The cost oscillates randomly.''], ''correctAnswer'': 2, ''explanation'': ''Lemma 3 states that during the course of the K-Means algorithm, the cost monotonically decreases. This means that with each iteration, the algorithm moves towards a lower cost (better clustering) until it converges.'', ''source'': {''chunkId'': 701, ''fileName'': ''Clustering.pdf'', ''chunkText'': '' the cost converge?\n\nLemma 3. During the course of the K-Means algorithm, the cost monotonically\ndecreases.\n\n(t) (t) (t) (t) Proof. Let z1 , ..., zK , C1 , ..., CK denote the centers and clusters at the start of\nthe tth iterate of K-Means. The first step o''}, ''grade'': None, ''feedback'': None}]}'
$ws.Cells.Item(2, 5).Value = 15
$ws.Cells.Item(2, 6).Value = 15
$ws.Cells.Item(2, 7).Value = 36.5
# Touch I2 so it exists as a blank cell (matches the added <c r="I2"/> in
# the target sheet) without picking up any new style/number format.
$ws.Cells.Item(2, 9).Borders.LineStyle = -4142

# --- Row 3: new full_course quiz row (quizId 58). ---
$ws.Cells.Item(3, 1).Value = "quiz_prompt_socratic"
$ws.Cells.Item(3, 2).Value = "full_course"
# Touch C3 so it exists as a blank cell (matches <c r="C3" t="inlineStr"/>).
$ws.Cells.Item(3, 3).Borders.LineStyle = -4142
$ws.Cells.Item(3, 4).Value = '{''quizId'': 58, ''courseId'': 11, ''materialIds'': [], ''section'': None, ''quizTitle'': ''Decision Tree Mastery Quiz'', ''questions'': [{''id'': 691, ''question'': ''Which package in R is used to implement the CART algorithm for decision trees?'', ''options'': [''e1071'', ''ada'', ''rpart'', ''randomForest''], ''correctAnswer'': 2, ''explanation'': ''The rpart package is specifically used for implementing the CART algorithm. e1071 contains SVM and Naive Bayes, and ada is for AdaBoost.'', ''source'': {''chunkId'': 138, ''fileName'': ''R for ML.pdf'', ''chunkText'': ''4.6 Na¨ıve Bayes\n\nInstall and load the e1071 package. Using the same notation as in Section 3.6, the command is:\n\n&gt; nB_model &lt;- naiveBayes(y ∼ x1 + x2, data=as.data.frame(cbind(y,x1,x2)))\n\n4.7 Decision Trees (CART)\n\nCART is implemented in the rpart package. Again using the formula, the command is\n\n&gt; cart_model &lt;- rpart(y ∼ x1 + x2, data=as.data.frame(cbind(y,x1,x2)), method="class")\n\nYou can use plot.rpart and text.rpart to plot the decision tree.\n\n4.8 AdaBoost\n\nThere are a number of different boosting func''}, ''grade'': None, ''feedback'': None}, {''id'': 692, ''question'': "What are the primary functions used to visualize decision trees created using the ''rpart'' package in R?", ''options'': [''plot() and image()'', ''plot.rpart() and text.rpart()'', ''draw() and label()'', ''visualize() and annotate()''], ''correctAnswer'': 1, ''explanation'': "The functions ''plot.rpart'' and ''text.rpart'' are specifically designed to plot and add text to decision trees generated by the ''rpart'' package. The other options are not valid functions for visualizing rpart objects.", ''source'': {''chunkId'': 138, ''fileName'': ''R for ML.pdf'', ''chunkText'': ''4.6 Na¨ıve Bayes\n\nInstall and load the e1071 package. Using the same notation as in Section 3.6, the command is:\n\n&gt; nB_model &lt;- naiveBayes(y ∼ x1 + x2, data=as.data.frame(cbind(y,x1,x2)))\n\n4.7 Decision Trees (CART)\n\nCART is implemented in the rpart package. Again using the formula, the command is\n\n&gt; cart_model &lt;- rpart(y ∼ x1 + x2, data=as.data.frame(cbind(y,x1,x2)), method="class")\n\nYou can use plot.rpart and text.rpart to plot the decision tree.\n\n4.8 AdaBoost\n\nThere are a number of different boosting func''}, ''grade'': None, ''feedback'': None}, {''id'': 693, ''question'': "In the context of decision tree pruning, what does a larger value of ''C'' (cost) typically imply?", ''options'': [''A more complex tree with better training accuracy'', ''A sparser tree with fewer leaves'', ''No impact on the tree structure'', ''A tree with higher misclassification error''], ''correctAnswer'': 1, ''explanation'': ''A larger value of C imposes a higher penalty for the number of leaves in the subtree, leading to a sparser, less complex tree. A smaller C would lead to better training accuracy but potentially overfitting.'', ''source'': {''chunkId'': 257, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''xity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. The second term is a regularization term. If we choose C to be large, the tree that minimizes the cost will be sparser. If C is small, the tree that minimizes the cost will have better training accuracy.\n\nX\n1[yi=leaf’s class] + C [#leaves in subtree] .''}, ''grade'': None, ''feedback'': None}, {''id'': 694, ''question'': ''Which of the following algorithms uses information gain for splitting nodes in a decision tree?'', ''options'': [''CART'', ''C4.5'', ''SVM'', ''Naive Bayes''], ''correctAnswer'': 1, ''explanation'': ''C4.5 uses information gain as its splitting criterion. CART uses the Gini index. SVM and Naive Bayes are not decision tree algorithms.'', ''source'': {''chunkId'': 245, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''Gin\n\nf p bet\n\nclassifica\n\ny,\n\nrop\n\ntion error for values o\n\nent\n\nof\n\nGra ph\n\ni i nde x, and mi s\n\nween 0 and 1\n\n.\n\nLearning, Springer, 2009.\n\nC4.5 uses information gain for splitting, and CART uses the Gini index. (CART only has binary splits.)\n\nPruning Let’s start with C4.5’s pruning. C4.5 recursively makes choices as to whether to prune on an attribute:\n\n• Option 1: leaving the tree as is\n\n• Option 2: replace that part of the tree with a leaf corresponding to the most frequent label in the data S going to''}, ''grade'': None, ''feedback'': None}, {''id'': 695, ''question'': ''What is the primary purpose of the regularization term in the cost function used for decision tree pruning?'', ''options'': [''To increase training accuracy'', ''To reduce the number of leaves in the tree'', ''To decrease misclassification error'', ''To increase the complexity of the tree''], ''correctAnswer'': 1, ''explanation'': ''The regularization term penalizes the number of leaves, thus encouraging simpler trees and preventing overfitting. It does not directly increase training accuracy or complexity.'', ''source'': {''chunkId'': 257, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''xity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. The second term is a regularization term. If we choose C to be large, the tree that minimizes the cost will be sparser. If C is small, the tree that minimizes the cost will have better training accuracy.\n\nX\n1[yi=leaf’s class] + C [#leaves in subtree] .''}, ''grade'': None, ''feedback'': None}, {''id'': 696, ''question'': "In CART''s regression trees, what value is assigned to f(x) in each leaf node?", ''options'': [''A linear function of x'', ''A constant value'', ''A probability distribution'', ''A complex non-linear function''], ''correctAnswer'': 1, ''explanation'': "In CART''s regression trees, each leaf node is assigned a constant value, representing the prediction for instances falling into that leaf. This contrasts with classification trees where leaves represent class labels.", ''source'': {''chunkId'': 260, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''CART’s Regression Trees\n\nFour-part figure in\n\ncl\n\nud ing partitions, a corresponding tree , and a perspective plot of the predicti on surface .\n\nLearning, Springer, 2009.\n\nCART decides which attributes to split and where to split them. In each leaf, we’re going to assign f(x) to be a constant.\n\nCan you guess what value to assign?\n\nConsider the empirical error, using the least squares loss:\n\nRtrain(f) =\nX (yi i −f(xi))2\n\n15''}, ''grade'': None, ''feedback'': None}, {''id'': 697, ''question'': ''What is the Gini index used for in the context of decision trees?'', ''options'': [''Measuring information gain'', ''Splitting nodes in CART'', ''Pruning the tree'', ''Visualizing the tree''], ''correctAnswer'': 1, ''explanation'': ''The Gini index is used by the CART algorithm to determine the best split at each node. Information gain is used by C4.5. Pruning and visualization are separate processes.'', ''source'': {''chunkId'': 245, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''Gin\n\nf p bet\n\nclassifica\n\ny,\n\nrop\n\ntion error for values o\n\nent\n\nof\n\nGra ph\n\ni i nde x, and mi s\n\nween 0 and 1\n\n.\n\nLearning, Springer, 2009.\n\nC4.5 uses information gain for splitting, and CART uses the Gini index. (CART only has binary splits.)\n\nPruning Let’s start with C4.5’s pruning. C4.5 recursively makes choices as to whether to prune on an attribute:\n\n• Option 1: leaving the tree as is\n\n• Option 2: replace that part of the tree with a leaf corresponding to the most frequent label in the data S going to''}, ''grade'': None, ''feedback'': None}, {''id'': 698, ''question'': ''Which of the following is NOT a typical application of decision trees?'', ''options'': [''Predicting house prices'', ''Automatic handwriting recognition'', ''Estimating probability of failure'', ''Finding patterns in large datasets (e.g., Diapers -&gt; Beer)''], ''correctAnswer'': 3, ''explanation'': ''Finding patterns in large datasets (like association rule mining) is typically done using algorithms like Apriori. Decision trees are well-suited for regression (predicting house prices), classification (handwriting recognition), and estimating probabilities.'', ''source'': {''chunkId'': 145, ''fileName'': ''Fundamentals of Learning.pdf'', ''chunkText'': ''• Input: {(xi, yi)}m\n\ni=1, xi ∈X , yi ∈ R\n\n• Output: f : X → R\n\n• predicting an individual’s income, predict house prices, predict stock prices, predict test scores\n\n5. Ranking (later) - in between classification and regression. Search engines\nuse ranking methods\n\n6. Density Estimation - predict conditional probabilities\n\ni=1, xi ∈X , yi ∈ {−1, 1}\n\n• {(xi, yi)}m\n\n• Output: f : X → [0, 1] as “close” to P (y = 1|x) as possible.\n\n• estimate probability of failure, probability to default on loan\n\nRule mining an''}, ''grade'': None, ''feedback'': None}, {''id'': 699, ''question'': "What is the purpose of holding out data in the process of choosing the cost complexity parameter ''C'' for pruning?", ''options'': [''To train the initial tree'', "To evaluate the performance of trees with different ''C'' values", ''To determine the splitting criteria'', ''To visualize the final tree''], ''correctAnswer'': 1, ''explanation'': "The holdout data is used to evaluate the performance of the pruned trees generated with different values of ''C''. This helps in selecting the ''C'' value that results in the best generalization performance.", ''source'': {''chunkId'': 258, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''imizes the cost will have better training accuracy.\n\nX\n1[yi=leaf’s class] + C [#leaves in subtree] .\n\nX\n\ncost(subtree) =\n\n̸\n\nves j xi∈\n\nlea\n\nleaf j\n\nWe could create a sequence of nested subtrees by gradually increasing C.\n\nDraw a picture\n\nNow we need to choose C. Here’s one way to do this:\n\n• Step 1: For each C, hold out some data, split, then prune, producing a tree for each C.\n\n• Step 2: see which tree is the best on the holdout data, choose C.\n\n• Step 3: use all data, use chosen C, run split, then prune ''}, ''grade'': None, ''feedback'': None}, {''id'': 700, ''question'': "What is the main difference between Option 1 and Option 2 in C4.5''s pruning process, as described in the material?", ''options'': [''Option 1 prunes the tree more aggressively than Option 2'', ''Option 1 leaves the tree as is, while Option 2 replaces a subtree with a leaf'', ''Option 1 uses information gain, while Option 2 uses the Gini index'', ''Option 1 is used for classification, while Option 2 is used for regression''], ''correctAnswer'': 1, ''explanation'': ''Option 1 in C4.5 pruning involves leaving the tree as it is, while Option 2 replaces a part of the tree with a leaf node corresponding to the most frequent label in the data going to that node. This is a key decision point in the pruning process.'', ''source'': {''chunkId'': 245, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''Gin\n\nf p bet\n\nclassifica\n\ny,\n\nrop\n\ntion error for values o\n\nent\n\nof\n\nGra ph\n\ni i nde x, and mi s\n\nween 0 and 1\n\n.\n\nLearning, Springer, 2009.\n\nC4.5 uses information gain for splitting, and CART uses the Gini index. (CART only has binary splits.)\n\nPruning Let’s start with C4.5’s pruning. C4.5 recursively makes choices as to whether to prune on an attribute:\n\n• Option 1: leaving the tree as is\n\n• Option 2: replace that part of the tree with a leaf corresponding to the most frequent label in the data S going to''}, ''grade'': None, ''feedback'': None}, {''id'': 701, ''question'': "What type of problem is predicting an individual''s income using a decision tree?", ''options'': [''Classification'', ''Regression'', ''Clustering'', ''Rule Mining''], ''correctAnswer'': 1, ''explanation'': ''Predicting a continuous value like income is a regression problem. Classification involves predicting a discrete category. Clustering groups data points, and rule mining discovers associations.'', ''source'': {''chunkId'': 145, ''fileName'': ''Fundamentals of Learning.pdf'', ''chunkText'': ''• Input: {(xi, yi)}m\n\ni=1, xi ∈X , yi ∈ R\n\n• Output: f : X → R\n\n• predicting an individual’s income, predict house prices, predict stock prices, predict test scores\n\n5. Ranking (later) - in between classification and regression. Search engines\nuse ranking methods\n\n6. Density Estimation - predict conditional probabilities\n\ni=1, xi ∈X , yi ∈ {−1, 1}\n\n• {(xi, yi)}m\n\n• Output: f : X → [0, 1] as “close” to P (y = 1|x) as possible.\n\n• estimate probability of failure, probability to default on loan\n\nRule mining an''}, ''grade'': None, ''feedback'': None}, {''id'': 702, ''question'': ''What is the role of cross-validation in the context of decision tree construction and pruning?'', ''options'': [''To determine the splitting criterion'', ''To visualize the tree'', ''To estimate the generalization performance for different pruning parameters'', ''To speed up the training process''], ''correctAnswer'': 2, ''explanation'': ''Cross-validation is used to estimate how well the tree will perform on unseen data for different values of the cost complexity parameter (C) during pruning. This helps in selecting the best C value to avoid overfitting.'', ''source'': {''chunkId'': 259, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''he best on the holdout data, choose C.\n\n• Step 3: use all data, use chosen C, run split, then prune to produce the final tree.\n\nThere are other ways to do this! You’ll use cross-validation in the homework.\n\nReview C4.5 and CART\n\n14''}, ''grade'': None, ''feedback'': None}, {''id'': 703, ''question'': ''Which of the following is a potential benefit of using decision trees in real-world applications?'', ''options'': [''They always outperform human experts'', ''They are always elegant and easy to understand'', ''They can replace complex, hand-designed rule systems'', ''They are not prone to overfitting''], ''correctAnswer'': 2, ''explanation'': "Decision trees can be used to replace complex rule systems, as demonstrated by BP''s GasOIL system. They don''t always outperform experts, aren''t always elegant, and can be prone to overfitting if not pruned properly.", ''source'': {''chunkId'': 223, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''e real examples (from Russell &amp; Norvig, Mitchell)\n\n• BP’s GasOIL system for separating gas and oil on offshore platforms - decision trees replaced a hand-designed rules system with 2500 rules. C4.5-based\nsystem outperformed human experts and saved BP millions. (1986)\n\n• learning to fly a Cessna on a flight simulator by watching human experts fly the simulator (1992)\n\n• can also learn to play tennis, analyze C-section risk, etc.\n\nHow to build a decision tree:\n\n• Start at the top of the tree.\n\n• Grow it by “s''}, ''grade'': None, ''feedback'': None}, {''id'': 704, ''question'': ''What is the main goal of the testing phase in machine learning, specifically in the context of decision trees?'', ''options'': [''To train the model'', ''To predict y for a new x and compare it to the true y'', ''To visualize the tree structure'', ''To select the best splitting attribute''], ''correctAnswer'': 1, ''explanation'': "The testing phase aims to evaluate the model''s ability to generalize to unseen data by predicting the output (y) for new input data (x) and comparing the prediction to the actual value of y.", ''source'': {''chunkId'': 147, ''fileName'': ''Fundamentals of Learning.pdf'', ''chunkText'': ''upervised learning\n\nTraining: training data are input, and model f is the output.\n\ni=1 =⇒ Algorithm =⇒ f.\n\n{(xi, yi)}m\n\nTesting: You want to predict y for a new x, where (x, y) comes from the same distribution as {(xi, yi)}m\n\ni=1.\n\nThat is, (x, y) ∼ D(X , Y) and each (xi, yi) ∼ D(X , Y).\n\nCompute f(x) and compare it to y. How well does f(x) match y? Measure goodness of f using a loss function R : Y × Y → R:\n\nRtest(f) = E(x,y)∼DR(f(x), y)\n\n=\nR(f(x), y)dD(x, y). (x,y)∼D\n\n2''}, ''grade'': None, ''feedback'': None}, {''id'': 705, ''question'': "What is the relationship between training data and the output ''f'' in supervised learning?", ''options'': ["''f'' is the input, and the training data is the output", "The training data is input, and model ''f'' is the output", "''f'' and the training data are independent of each other", "''f'' is only used for visualization"], ''correctAnswer'': 1, ''explanation'': "In supervised learning, the training data ({(xi, yi)}m i=1) is used as input to an algorithm, and the resulting model or function ''f'' is the output. This ''f'' is then used to make predictions on new data.", ''source'': {''chunkId'': 147, ''fileName'': ''Fundamentals of Learning.pdf'', ''chunkText'': ''upervised learning\n\nTraining: training data are input, and model f is the output.\n\ni=1 =⇒ Algorithm =⇒ f.\n\n{(xi, yi)}m\n\nTesting: You want to predict y for a new x, where (x, y) comes from the same distribution as {(xi, yi)}m\n\ni=1.\n\nThat is, (x, y) ∼ D(X , Y) and each (xi, yi) ∼ D(X , Y).\n\nCompute f(x) and compare it to y. How well does f(x) match y? Measure goodness of f using a loss function R : Y × Y → R:\n\nRtest(f) = E(x,y)∼DR(f(x), y)\n\n=\nR(f(x), y)dD(x, y). (x,y)∼D\n\n2''}, ''grade'': None, ''feedback'': None}]}'
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 15
$ws.Cells.Item(3, 7).Value = 31.89
# Touch H3 / I3 so they exist as blank cells.
$ws.Cells.Item(3, 8).Borders.LineStyle = -4142
$ws.Cells.Item(3, 9).Borders.LineStyle = -4142

# --- Row 4: new topic quiz row (quizId 59), with topic = "decision trees". ---
$ws.Cells.Item(4, 1).Value = "quiz_prompt_socratic"
$ws.Cells.Item(4, 2).Value = "topic"
# Touch C4 so it exists as a blank cell.
$ws.Cells.Item(4, 3).Borders.LineStyle = -4142
$ws.Cells.Item(4, 4).Value = '{''quizId'': 59, ''courseId'': 11, ''materialIds'': [], ''section'': None, ''quizTitle'': ''Decision Tree Mastery Quiz'', ''questions'': [{''id'': 706, ''question'': "What is the primary goal of the CART algorithm''s minimal cost complexity pruning?", ''options'': [''To minimize misclassification errors on the training data at all costs.'', ''To find the smallest possible tree, regardless of accuracy.'', ''To balance misclassification error with tree complexity, preventing overfitting.'', ''To maximize the number of leaves in the tree for better interpretability.''], ''correctAnswer'': 2, ''explanation'': "CART''s minimal cost complexity pruning aims to find a balance between minimizing misclassification errors and keeping the tree simple. This helps to prevent overfitting by penalizing complex trees with many leaves. A large C value favors sparser trees, while a small C value favors better training accuracy. The other options represent extremes that are not the goal of the algorithm.", ''source'': {''chunkId'': 256, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''CART - Classification and Regression Trees (Breiman, Freedman, Olshen, Stone,\n1984)\n\nDoes only binary splits, not multiway splits (less interpretable, but simplifies splitting criteria).\n\nLet’s do classification first, and regression later. For splitting, CART uses the Gini index. The Gini index again is\n\nvariance of Bin(n, p)\np(1 −p) =\n= variance of Bernoulli(p).\nn For pruning, CART uses “minimal cost complexity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. Th''}, ''grade'': None, ''feedback'': None}, {''id'': 707, ''question'': "In the context of decision trees, what does ''Information Gain'' represent?", ''options'': [''The amount of data added to the tree after a split.'', ''The expected reduction in entropy due to branching on an attribute.'', ''The increase in entropy after a split.'', ''The number of positive examples in a dataset.''], ''correctAnswer'': 1, ''explanation'': ''Information Gain quantifies the reduction in entropy (uncertainty) achieved by splitting a dataset on a particular attribute. A higher information gain indicates a more effective split. The other options are incorrect because they either misrepresent the concept of entropy or describe unrelated aspects of decision trees.'', ''source'': {''chunkId'': 236, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''gatives.\n\nThe training probabilities\n\nin branch j are:\n\n#posj #negj ,\n#posj + #negj #posj + #negj\n\n.\n\nThe Information Gain is calculated like this:\n\nGain(S, A) = expected reduction in entropy due to branching on attribute A\n= original\n\nentropy −entropy after branching #pos #neg\n= H\n,\n#pos + #neg #pos + #neg\n\nJ\n\n−\nX #posj + #negj\n#posj #negj H , .''}, ''grade'': None, ''feedback'': None}, {''id'': 708, ''question'': ''How does CART handle splitting continuous attributes?'', ''options'': [''It uses multiway splits to divide the attribute into multiple ranges.'', ''It discretizes the attribute into predefined bins before splitting.'', ''It performs binary splits by finding the optimal split point through a line search.'', ''It ignores continuous attributes and only uses categorical attributes.''], ''correctAnswer'': 2, ''explanation'': "CART handles continuous attributes by performing binary splits. It searches for the optimal split point ''s'' that minimizes the impurity (e.g., Gini index) of the resulting subsets. This involves a line search over possible values of ''s''. Multiway splits are not used in CART, and continuous attributes are not ignored.", ''source'': {''chunkId'': 263, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''dily want attribute A and split point s solving the following.\n\n\uf8f9\n\n\uf8ee\n\n\uf8f0min (yi\nC1)2 + min\n(yi C2)2 . A, s C1 x\n\nX\n\nX\n\nC2 leaf| (A) − −\n\nmin\n\n\uf8fb\n\n( i leaf for each attribute A ∈{\n|x A)&gt;s\ndo }\n\ni∈{ x ≤s} x\n\na\n\nlinesearch over s\n\nThe first term means that we’ll choose the optimal C1 = y¯{leaf|x(A)≤s . The second\n}\nterm means we’ll choose C2 = y¯{leaf x(A)&gt;s .\n| }\n\nFor pruning, again CART does minimal cost complexity pruning:\n\ni∈S\n−y¯Sj)2 + C[# leaves in tree]\n\nX\n\nX yi\n\ncost =\n(\n\nves j x\n\nlea\n\nj\n\n16''}, ''grade'': None, ''feedback'': None}, {''id'': 709, ''question'': ''What is the Gini index used for in the CART algorithm?'', ''options'': [''To measure the impurity of a node or a split.'', ''To determine the depth of the tree.'', ''To calculate the information gain of an attribute.'', ''To prune the tree after it has been built.''], ''correctAnswer'': 0, ''explanation'': ''The Gini index is a measure of impurity used in CART for classification tasks. It quantifies the probability of misclassifying a randomly chosen element in a dataset if it were randomly labeled according to the distribution of labels in the subset. The other options are incorrect because they describe different aspects or steps in the decision tree process.'', ''source'': {''chunkId'': 256, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''CART - Classification and Regression Trees (Breiman, Freedman, Olshen, Stone,\n1984)\n\nDoes only binary splits, not multiway splits (less interpretable, but simplifies splitting criteria).\n\nLet’s do classification first, and regression later. For splitting, CART uses the Gini index. The Gini index again is\n\nvariance of Bin(n, p)\np(1 −p) =\n= variance of Bernoulli(p).\nn For pruning, CART uses “minimal cost complexity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. Th''}, ''grade'': None, ''feedback'': None}, {''id'': 710, ''question'': ''What is a key difference between C4.5 and CART in terms of splitting?'', ''options'': [''C4.5 only handles binary splits, while CART handles multiway splits.'', ''CART only handles binary splits, while C4.5 can handle multiway splits.'', ''C4.5 uses the Gini index, while CART uses entropy for splitting.'', ''CART is more elegant and interpretable than C4.5.''], ''correctAnswer'': 1, ''explanation'': ''CART is characterized by its use of only binary splits, while C4.5 can handle multiway splits. This is a fundamental difference in how the two algorithms construct decision trees. While the material suggests C4.5 and CART are not elegant, it does not suggest that CART is more elegant or interpretable than C4.5.'', ''source'': {''chunkId'': 256, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''CART - Classification and Regression Trees (Breiman, Freedman, Olshen, Stone,\n1984)\n\nDoes only binary splits, not multiway splits (less interpretable, but simplifies splitting criteria).\n\nLet’s do classification first, and regression later. For splitting, CART uses the Gini index. The Gini index again is\n\nvariance of Bin(n, p)\np(1 −p) =\n= variance of Bernoulli(p).\nn For pruning, CART uses “minimal cost complexity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. Th''}, ''grade'': None, ''feedback'': None}, {''id'': 711, ''question'': ''Which of the following is NOT a typical application of decision trees?'', ''options'': [''Predicting house prices.'', ''Analyzing C-section risk.'', ''Classifying program behavior as normal or intrusive.'', ''Estimating probability of failure.''], ''correctAnswer'': 2, ''explanation'': ''Decision trees can be used for classification and regression. Predicting house prices is a regression task. Analyzing C-section risk and estimating probability of failure are classification tasks. Classifying program behavior as normal or intrusive is typically done using K-Nearest Neighbors.'', ''source'': {''chunkId'': 223, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''e real examples (from Russell &amp; Norvig, Mitchell)\n\n• BP’s GasOIL system for separating gas and oil on offshore platforms - decision trees replaced a hand-designed rules system with 2500 rules. C4.5-based\nsystem outperformed human experts and saved BP millions. (1986)\n\n• learning to fly a Cessna on a flight simulator by watching human experts fly the simulator (1992)\n\n• can also learn to play tennis, analyze C-section risk, etc.\n\nHow to build a decision tree:\n\n• Start at the top of the tree.\n\n• Grow it by “s''}, ''grade'': None, ''feedback'': None}, {''id'': 712, ''question'': "What is the purpose of setting a higher value for ''C'' in CART''s cost complexity pruning?", ''options'': [''To achieve better training accuracy.'', ''To obtain a sparser (smaller) tree.'', ''To minimize misclassification errors on the training data.'', ''To increase the number of leaves in the tree.''], ''correctAnswer'': 1, ''explanation'': "In CART''s cost complexity pruning, ''C'' is a regularization parameter. A higher value of ''C'' penalizes tree complexity more heavily, leading to a sparser (smaller) tree. Conversely, a smaller ''C'' allows for more complex trees that may achieve better training accuracy but risk overfitting.", ''source'': {''chunkId'': 257, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''xity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. The second term is a regularization term. If we choose C to be large, the tree that minimizes the cost will be sparser. If C is small, the tree that minimizes the cost will have better training accuracy.\n\nX\n1[yi=leaf’s class] + C [#leaves in subtree] .''}, ''grade'': None, ''feedback'': None}, {''id'': 713, ''question'': ''Which of the following is a disadvantage of using decision trees?'', ''options'': [''They are difficult to interpret.'', ''They cannot handle categorical data.'', ''They are prone to overfitting if not pruned.'', ''They can only be used for regression tasks.''], ''correctAnswer'': 3, ''explanation'': ''Decision trees, if left unpruned, can become overly complex and memorize the training data, leading to poor performance on unseen data (overfitting). They are generally easy to interpret, can handle categorical data, and can be used for both classification and regression tasks.'', ''source'': {''chunkId'': 257, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''xity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. The second term is a regularization term. If we choose C to be large, the tree that minimizes the cost will be sparser. If C is small, the tree that minimizes the cost will have better training accuracy.\n\nX\n1[yi=leaf’s class] + C [#leaves in subtree] .''}, ''grade'': None, ''feedback'': None}, {''id'': 714, ''question'': "What is the role of a ''holdout set'' in the context of decision tree construction?", ''options'': [''To train the decision tree model.'', "To select the best value for the complexity parameter ''C'' during pruning.", ''To estimate the information gain of different attributes.'', ''To visualize the decision tree.''], ''correctAnswer'': 1, ''explanation'': "A holdout set (or validation set) is used to evaluate the performance of the decision tree with different values of the complexity parameter ''C'' during pruning. The value of ''C'' that results in the best performance on the holdout set is then chosen for the final model. The training data is used to train the model, and the information gain is calculated during the tree building process.", ''source'': {''chunkId'': 259, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''he best on the holdout data, choose C.\n\n• Step 3: use all data, use chosen C, run split, then prune to produce the final tree.\n\nThere are other ways to do this! You’ll use cross-validation in the homework.\n\nReview C4.5 and CART\n\n14''}, ''grade'': None, ''feedback'': None}, {''id'': 715, ''question'': ''What is the purpose of calculating upper bounds on the probability of error during decision tree pruning (as mentioned with C4.5)?'', ''options'': [''To estimate the accuracy of the unpruned tree.'', ''To compare different pruning options (e.g., replacing a subtree with a leaf or a branch).'', ''To determine the optimal depth of the tree.'', ''To calculate the information gain of each attribute.''], ''correctAnswer'': 1, ''explanation'': ''Calculating upper bounds on the probability of error for different pruning options (e.g., replacing a subtree with a leaf or a branch) allows the algorithm to compare these options and choose the one with the lowest estimated error. This helps to prevent overfitting and improve the generalization performance of the tree. The other options describe different aspects of decision tree construction or evaluation.'', ''source'': {''chunkId'': 246, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''ce that part of the tree with a leaf corresponding to the most frequent label in the data S going to that part of the tree.\n\n• Option 3: replace that part of the tree with one of its subtrees, corresponding to the most common branch in the split\n\nDemo To figure out which decision to make, C4.5 computes upper bounds on the probability of error for each option. I’ll show you how to do that shortly.\n\n• Prob of error for Option 1 ≤UpperBound1\n\n• Prob of error for Option 2 ≤UpperBound2\n\n10''}, ''grade'': None, ''feedback'': None}, {''id'': 716, ''question'': "Which of the following best describes the ''minimal cost complexity'' pruning method used by CART?", ''options'': [''A method that only considers the misclassification error when pruning.'', ''A method that balances misclassification error with the complexity of the tree.'', ''A method that always chooses the smallest possible tree.'', ''A method that prioritizes interpretability over accuracy.''], ''correctAnswer'': 1, ''explanation'': "Minimal cost complexity pruning in CART aims to find the optimal balance between minimizing misclassification error and keeping the tree simple. It adds a penalty term for the number of leaves in the tree, controlled by the complexity parameter ''C''. This prevents overfitting by favoring simpler trees that generalize better to unseen data. The other options are incorrect because they represent extremes or misinterpret the goal of the method.", ''source'': {''chunkId'': 256, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''CART - Classification and Regression Trees (Breiman, Freedman, Olshen, Stone,\n1984)\n\nDoes only binary splits, not multiway splits (less interpretable, but simplifies splitting criteria).\n\nLet’s do classification first, and regression later. For splitting, CART uses the Gini index. The Gini index again is\n\nvariance of Bin(n, p)\np(1 −p) =\n= variance of Bernoulli(p).\nn For pruning, CART uses “minimal cost complexity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. Th''}, ''grade'': None, ''feedback'': None}, {''id'': 717, ''question'': ''In the context of decision trees, what is the purpose of cross-validation?'', ''options'': [''To speed up the training process.'', ''To estimate the generalization performance of the model and select the best hyperparameters.'', ''To visualize the decision tree.'', ''To simplify the decision tree.''], ''correctAnswer'': 1, ''explanation'': "Cross-validation is a technique used to estimate how well a model will generalize to unseen data. It involves splitting the data into multiple folds, training the model on some folds, and evaluating it on the remaining fold. This process is repeated multiple times, and the results are averaged to obtain a more robust estimate of the model''s performance. It can also be used to select the best hyperparameters, such as the complexity parameter ''C'' in CART.", ''source'': {''chunkId'': 259, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''he best on the holdout data, choose C.\n\n• Step 3: use all data, use chosen C, run split, then prune to produce the final tree.\n\nThere are other ways to do this! You’ll use cross-validation in the homework.\n\nReview C4.5 and CART\n\n14''}, ''grade'': None, ''feedback'': None}, {''id'': 718, ''question'': ''Which of the following real-world applications directly benefited from the use of decision trees, as mentioned in the provided material?'', ''options'': [''Spam email filtering.'', ''Fraud detection in credit card transactions.'', ''Separating gas and oil on offshore platforms.'', ''Automatic translation of languages.''], ''correctAnswer'': 2, ''explanation'': "The material mentions that BP''s GasOIL system for separating gas and oil on offshore platforms used decision trees (specifically, a C4.5-based system) to replace a hand-designed rules system, resulting in significant cost savings. The other options are not mentioned in the context of decision tree applications.", ''source'': {''chunkId'': 223, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''e real examples (from Russell &amp; Norvig, Mitchell)\n\n• BP’s GasOIL system for separating gas and oil on offshore platforms - decision trees replaced a hand-designed rules system with 2500 rules. C4.5-based\nsystem outperformed human experts and saved BP millions. (1986)\n\n• learning to fly a Cessna on a flight simulator by watching human experts fly the simulator (1992)\n\n• can also learn to play tennis, analyze C-section risk, etc.\n\nHow to build a decision tree:\n\n• Start at the top of the tree.\n\n• Grow it by “s''}, ''grade'': None, ''feedback'': None}, {''id'': 719, ''question'': ''What is the main advantage of using binary splits in CART?'', ''options'': [''They are more interpretable than multiway splits.'', ''They can handle categorical attributes with many values more efficiently.'', ''They simplify the splitting criteria.'', ''They always result in smaller trees.''], ''correctAnswer'': 2, ''explanation'': ''The material states that binary splits in CART simplify the splitting criteria. While binary splits can sometimes be more interpretable, the primary reason for using them is to make the splitting process more manageable, especially when dealing with continuous or categorical attributes with many values. They do not necessarily result in smaller trees.'', ''source'': {''chunkId'': 256, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''CART - Classification and Regression Trees (Breiman, Freedman, Olshen, Stone,\n1984)\n\nDoes only binary splits, not multiway splits (less interpretable, but simplifies splitting criteria).\n\nLet’s do classification first, and regression later. For splitting, CART uses the Gini index. The Gini index again is\n\nvariance of Bin(n, p)\np(1 −p) =\n= variance of Bernoulli(p).\nn For pruning, CART uses “minimal cost complexity.”\n\nEach subtree is assigned a cost. The first term in the cost is a misclassification error. Th''}, ''grade'': None, ''feedback'': None}, {''id'': 720, ''question'': ''What is the relationship between entropy and information gain?'', ''options'': [''Information gain is the same as entropy.'', ''Information gain is the inverse of entropy.'', ''Information gain is the expected reduction in entropy due to splitting on an attribute.'', ''Entropy is the expected reduction in information gain due to splitting on an attribute.''], ''correctAnswer'': 2, ''explanation'': ''Information gain measures how much the entropy (uncertainty) of a dataset is reduced after splitting on a particular attribute. A higher information gain indicates a more effective split in terms of reducing uncertainty. The other options misrepresent the relationship between entropy and information gain.'', ''source'': {''chunkId'': 236, ''fileName'': ''Decision trees.pdf'', ''chunkText'': ''gatives.\n\nThe training probabilities\n\nin branch j are:\n\n#posj #negj ,\n#posj + #negj #posj + #negj\n\n.\n\nThe Information Gain is calculated like this:\n\nGain(S, A) = expected reduction in entropy due to branching on attribute A\n= original\n\nentropy −entropy after branching #pos #neg\n= H\n,\n#pos + #neg #pos + #neg\n\nJ\n\n−\nX #posj + #negj\n#posj #negj H , .''}, ''grade'': None, ''feedback'': None}]}'
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 15
$ws.Cells.Item(4, 7).Value = 39.38
# Touch H4 so it exists as a blank cell.
$ws.Cells.Item(4, 8).Borders.LineStyle = -4142
$ws.Cells.Item(4, 9).Value = "decision trees"
